$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell 'D2' '42.779.12'
$ws.Range('E2').Value = '  +0.28%  '
Set-TextCell 'D3' '2.558.19'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  -0.04%  '
Set-TextCell 'D5' '310.01'
$ws.Range('E5').Value = '  -1.91%  '
Set-TextCell 'D6' '98.44'
$ws.Range('E6').Value = '  +2.51%  '
$ws.Range('E7').Value = '  -0.99%  '
$ws.Range('E8').Value = '  -0.06%  '
Set-TextCell 'D9' '0.530'
$ws.Range('E9').Value = '  -0.31%  '
Set-TextCell 'D10' '35.68'
$ws.Range('E10').Value = '  -0.04%  '
Set-TextCell 'D11' '0.0807'
$ws.Range('E11').Value = '  +0.18%  '
Set-TextCell 'D12' '7.43'
$ws.Range('E12').Value = '  -1.24%  '
Set-TextCell 'D13' '2.958.54'
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('E15').Value = '  +4.91%  '
Set-TextCell 'D16' '2.562.47'
$ws.Range('E16').Value = '  +2.25%  '
Set-TextCell 'D17' '0.838'
$ws.Range('E17').Value = '  -1.31%  '
Set-TextCell 'D18' '42.797.88'
$ws.Range('E18').Value = '  +0.02%  '
Set-TextCell 'D19' '6.74'
$ws.Range('E19').Value = '  -0.82%  '
Set-TextCell 'D20' '0.0₃0959'
$ws.Range('E20').Value = '  -0.38%  '
Set-TextCell 'D21' '12.36'
$ws.Range('E21').Value = '  -3.23%  '
Set-TextCell 'D22' '69.32'
$ws.Range('E22').Value = '  -0.33%  '
Set-TextCell 'D23' '247.51'
$ws.Range('E23').Value = '  -1.36%  '
Set-TextCell 'D24' '2.92'
$ws.Range('E24').Value = '  -0.90%  '
Set-TextCell 'D25' '2.04'
$ws.Range('E25').Value = '  -0.59%  '
Set-TextCell 'D26' '26.87'
$ws.Range('E26').Value = '  +1.77%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  -0.22%  '
Set-TextCell 'D29' '39.89'
$ws.Range('E29').Value = '  -2.22%  '
Set-TextCell 'D30' '10.16'
$ws.Range('E30').Value = '  -2.24%  '
Set-TextCell 'D31' '158.36'
$ws.Range('E31').Value = '  +0.60%  '
Set-TextCell 'D32' '5.76'
$ws.Range('E32').Value = '  -2.53%  '
Set-TextCell 'D33' '0.0796'
$ws.Range('E33').Value = '  +1.79%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 'D34' '2.67'
$ws.Range('E34').Value = '  -1.82%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D35' '2.09'
$ws.Range('E35').Value = '  -2.68%  '
$ws.Range('E36').Value = '  -1.57%  '
Set-TextCell 'D37' '18.78'
$ws.Range('E37').Value = '  -0.43%  '
$ws.Range('E38').Value = '  +12.66%  '
Set-TextCell 'D39' '0.111'
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('E40').Value = '  -0.50%  '
Set-TextCell 'D41' '23.06'
$ws.Range('E41').Value = '  +2.55%  '
Set-TextCell 'D42' '4.08'
$ws.Range('E42').Value = '  +6.87%  '
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('E44').Value = '  -1.06%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 'D45' '1.991.70'
$ws.Range('E45').Value = '  -2.05%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D46' '3.20'
$ws.Range('E46').Value = '  -1.91%  '
Set-TextCell 'D47' '8.99'
$ws.Range('E47').Value = '  -0.84%  '
Set-TextCell 'D48' '2.810.56'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('B49').Value = 'BitcoinSV'
$ws.Range('C49').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextCell 'D49' '81.20'
$ws.Range('E49').Value = '  -3.62%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D50' '0.193'
$ws.Range('E50').Value = '  +1.76%  '
Set-TextCell 'D51' '73.52'
$ws.Range('E51').Value = '  -2.56%  '
